{"js": "// Update the date line and all two-digit multiplication problem/answer cells.\n// Each original cell value is unique in the document, so a direct\n// search-and-replace (matchCase to avoid accidental partial matches) is safe.\nconst replacements = [\n  [\"2024-07-05 Friday\", \"2024-07-06 Saturday\"],\n  [\"24\u00d785=2040\", \"63\u00d775=4725\"],\n  [\"13\u00d783=1079\", \"58\u00d753=3074\"],\n  [\"95\u00d733=3135\", \"77\u00d715=1155\"],\n  [\"89\u00d781=7209\", \"47\u00d795=4465\"],\n  [\"83\u00d793=7719\", \"89\u00d793=8277\"],\n  [\"63\u00d737=2331\", \"48\u00d722=1056\"],\n  [\"56\u00d738=2128\", \"11\u00d746=506\"],\n  [\"95\u00d794=8930\", \"20\u00d770=1400\"],\n  [\"25\u00d797=2425\", \"91\u00d721=1911\"],\n  [\"86\u00d783=7138\", \"25\u00d764=1600\"],\n  [\"98\u00d795=9310\", \"96\u00d717=1632\"],\n  [\"26\u00d797=2522\", \"35\u00d763=2205\"],\n  [\"63\u00d760=3780\", \"99\u00d718=1782\"],\n  [\"11\u00d724=264\", \"56\u00d724=1344\"],\n  [\"47\u00d721=987\", \"41\u00d784=3444\"],\n  [\"65\u00d726=1690\", \"23\u00d750=1150\"],\n  [\"41\u00d720=820\", \"94\u00d784=7896\"],\n  [\"93\u00d724=2232\", \"62\u00d795=5890\"],\n  [\"63\u00d726=1638\", \"12\u00d787=1044\"],\n  [\"79\u00d760=4740\", \"99\u00d711=1089\"],\n  [\"32\u00d713=416\", \"78\u00d727=2106\"],\n  [\"39\u00d714=546\", \"21\u00d739=819\"],\n  [\"64\u00d779=5056\", \"28\u00d777=2156\"],\n  [\"56\u00d719=1064\", \"17\u00d742=714\"],\n  [\"54\u00d789=4806\", \"62\u00d788=5456\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, 'Replace');\n  }\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n}\n", "ps1": "# Update the date line and all two-digit multiplication problem/answer cells.\n# Each original cell value is unique in the document, so a simple\n# Find/Replace (MatchCase, no wildcards) for each pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-07-05 Friday\", \"2024-07-06 Saturday\")\n    ,@(\"24\u00d785=2040\", \"63\u00d775=4725\")\n    ,@(\"13\u00d783=1079\", \"58\u00d753=3074\")\n    ,@(\"95\u00d733=3135\", \"77\u00d715=1155\")\n    ,@(\"89\u00d781=7209\", \"47\u00d795=4465\")\n    ,@(\"83\u00d793=7719\", \"89\u00d793=8277\")\n    ,@(\"63\u00d737=2331\", \"48\u00d722=1056\")\n    ,@(\"56\u00d738=2128\", \"11\u00d746=506\")\n    ,@(\"95\u00d794=8930\", \"20\u00d770=1400\")\n    ,@(\"25\u00d797=2425\", \"91\u00d721=1911\")\n    ,@(\"86\u00d783=7138\", \"25\u00d764=1600\")\n    ,@(\"98\u00d795=9310\", \"96\u00d717=1632\")\n    ,@(\"26\u00d797=2522\", \"35\u00d763=2205\")\n    ,@(\"63\u00d760=3780\", \"99\u00d718=1782\")\n    ,@(\"11\u00d724=264\", \"56\u00d724=1344\")\n    ,@(\"47\u00d721=987\", \"41\u00d784=3444\")\n    ,@(\"65\u00d726=1690\", \"23\u00d750=1150\")\n    ,@(\"41\u00d720=820\", \"94\u00d784=7896\")\n    ,@(\"93\u00d724=2232\", \"62\u00d795=5890\")\n    ,@(\"63\u00d726=1638\", \"12\u00d787=1044\")\n    ,@(\"79\u00d760=4740\", \"99\u00d711=1089\")\n    ,@(\"32\u00d713=416\", \"78\u00d727=2106\")\n    ,@(\"39\u00d714=546\", \"21\u00d739=819\")\n    ,@(\"64\u00d779=5056\", \"28\u00d777=2156\")\n    ,@(\"56\u00d719=1064\", \"17\u00d742=714\")\n    ,@(\"54\u00d789=4806\", \"62\u00d788=5456\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceAll\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
